# Add rows 157-179 to the "pelada" sheet (new players' match stats for this week)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(157, 1).Value = "Jorge"
$ws.Cells.Item(157, 3).Value = 3
$ws.Cells.Item(157, 4).Value = 3
$ws.Cells.Item(157, 5).Value = 3
$ws.Cells.Item(157, 6).Value = 1
$ws.Cells.Item(157, 7).Value = 1
$ws.Cells.Item(157, 8).Value = 0
$ws.Cells.Item(157, 9).Value = 0
$ws.Cells.Item(157, 10).Value = 0
$ws.Cells.Item(157, 11).Value = 0
$ws.Cells.Item(157, 12).Value = 0

$ws.Cells.Item(158, 1).Value = "Juscielio"
$ws.Cells.Item(158, 3).Value = 3
$ws.Cells.Item(158, 4).Value = 3
$ws.Cells.Item(158, 5).Value = 3
$ws.Cells.Item(158, 6).Value = 2
$ws.Cells.Item(158, 7).Value = 1
$ws.Cells.Item(158, 8).Value = 0
$ws.Cells.Item(158, 9).Value = 0
$ws.Cells.Item(158, 10).Value = 0
$ws.Cells.Item(158, 11).Value = 0
$ws.Cells.Item(158, 12).Value = 0

$ws.Cells.Item(159, 1).Value = "Eduardo"
$ws.Cells.Item(159, 3).Value = 3
$ws.Cells.Item(159, 4).Value = 3
$ws.Cells.Item(159, 5).Value = 3
$ws.Cells.Item(159, 6).Value = 3
$ws.Cells.Item(159, 7).Value = 1
$ws.Cells.Item(159, 8).Value = 0
$ws.Cells.Item(159, 9).Value = 0
$ws.Cells.Item(159, 10).Value = 0
$ws.Cells.Item(159, 11).Value = 0
$ws.Cells.Item(159, 12).Value = 0

$ws.Cells.Item(160, 1).Value = "Heider"
$ws.Cells.Item(160, 3).Value = 3
$ws.Cells.Item(160, 4).Value = 3
$ws.Cells.Item(160, 5).Value = 3
$ws.Cells.Item(160, 6).Value = 0
$ws.Cells.Item(160, 7).Value = 1
$ws.Cells.Item(160, 8).Value = 0
$ws.Cells.Item(160, 9).Value = 0
$ws.Cells.Item(160, 10).Value = 0
$ws.Cells.Item(160, 11).Value = 0
$ws.Cells.Item(160, 12).Value = 0

$ws.Cells.Item(161, 1).Value = "Marcelão"
$ws.Cells.Item(161, 3).Value = 3
$ws.Cells.Item(161, 4).Value = 3
$ws.Cells.Item(161, 5).Value = 3
$ws.Cells.Item(161, 6).Value = 1
$ws.Cells.Item(161, 7).Value = 1
$ws.Cells.Item(161, 8).Value = 0
$ws.Cells.Item(161, 9).Value = 0
$ws.Cells.Item(161, 10).Value = 0
$ws.Cells.Item(161, 11).Value = 0
$ws.Cells.Item(161, 12).Value = 0

$ws.Cells.Item(162, 1).Value = "Caio"
$ws.Cells.Item(162, 3).Value = 1
$ws.Cells.Item(162, 4).Value = 4
$ws.Cells.Item(162, 5).Value = 3
$ws.Cells.Item(162, 6).Value = 1
$ws.Cells.Item(162, 7).Value = 1
$ws.Cells.Item(162, 8).Value = 0
$ws.Cells.Item(162, 9).Value = 1
$ws.Cells.Item(162, 10).Value = 0
$ws.Cells.Item(162, 11).Value = 0
$ws.Cells.Item(162, 12).Value = 0

$ws.Cells.Item(163, 1).Value = "Rulli"
$ws.Cells.Item(163, 3).Value = 1
$ws.Cells.Item(163, 4).Value = 4
$ws.Cells.Item(163, 5).Value = 3
$ws.Cells.Item(163, 6).Value = 0
$ws.Cells.Item(163, 7).Value = 1
$ws.Cells.Item(163, 8).Value = 0
$ws.Cells.Item(163, 9).Value = 1
$ws.Cells.Item(163, 10).Value = 0
$ws.Cells.Item(163, 11).Value = 0
$ws.Cells.Item(163, 12).Value = 0

$ws.Cells.Item(164, 1).Value = "Michel"
$ws.Cells.Item(164, 3).Value = 1
$ws.Cells.Item(164, 4).Value = 4
$ws.Cells.Item(164, 5).Value = 3
$ws.Cells.Item(164, 6).Value = 0
$ws.Cells.Item(164, 7).Value = 1
$ws.Cells.Item(164, 8).Value = 0
$ws.Cells.Item(164, 9).Value = 1
$ws.Cells.Item(164, 10).Value = 0
$ws.Cells.Item(164, 11).Value = 0
$ws.Cells.Item(164, 12).Value = 0

$ws.Cells.Item(165, 1).Value = "Corinthiano"
$ws.Cells.Item(165, 3).Value = 1
$ws.Cells.Item(165, 4).Value = 4
$ws.Cells.Item(165, 5).Value = 3
$ws.Cells.Item(165, 6).Value = 2
$ws.Cells.Item(165, 7).Value = 1
$ws.Cells.Item(165, 8).Value = 0
$ws.Cells.Item(165, 9).Value = 1
$ws.Cells.Item(165, 10).Value = 0
$ws.Cells.Item(165, 11).Value = 0
$ws.Cells.Item(165, 12).Value = 0

$ws.Cells.Item(166, 1).Value = "Miqueias"
$ws.Cells.Item(166, 3).Value = 1
$ws.Cells.Item(166, 4).Value = 4
$ws.Cells.Item(166, 5).Value = 3
$ws.Cells.Item(166, 6).Value = 2
$ws.Cells.Item(166, 7).Value = 1
$ws.Cells.Item(166, 8).Value = 0
$ws.Cells.Item(166, 9).Value = 1
$ws.Cells.Item(166, 10).Value = 0
$ws.Cells.Item(166, 11).Value = 0
$ws.Cells.Item(166, 12).Value = 0

$ws.Cells.Item(167, 1).Value = "Nenzinho"
$ws.Cells.Item(167, 3).Value = 4
$ws.Cells.Item(167, 4).Value = 2
$ws.Cells.Item(167, 5).Value = 2
$ws.Cells.Item(167, 6).Value = 0
$ws.Cells.Item(167, 7).Value = 1
$ws.Cells.Item(167, 8).Value = 1
$ws.Cells.Item(167, 9).Value = 0
$ws.Cells.Item(167, 10).Value = 0
$ws.Cells.Item(167, 11).Value = 0
$ws.Cells.Item(167, 12).Value = 0

$ws.Cells.Item(168, 1).Value = "Fabinho"
$ws.Cells.Item(168, 3).Value = 4
$ws.Cells.Item(168, 4).Value = 2
$ws.Cells.Item(168, 5).Value = 2
$ws.Cells.Item(168, 6).Value = 2
$ws.Cells.Item(168, 7).Value = 1
$ws.Cells.Item(168, 8).Value = 1
$ws.Cells.Item(168, 9).Value = 0
$ws.Cells.Item(168, 10).Value = 0
$ws.Cells.Item(168, 11).Value = 0
$ws.Cells.Item(168, 12).Value = 0

$ws.Cells.Item(169, 1).Value = "Cabeleira"
$ws.Cells.Item(169, 3).Value = 4
$ws.Cells.Item(169, 4).Value = 2
$ws.Cells.Item(169, 5).Value = 2
$ws.Cells.Item(169, 6).Value = 4
$ws.Cells.Item(169, 7).Value = 1
$ws.Cells.Item(169, 8).Value = 1
$ws.Cells.Item(169, 9).Value = 0
$ws.Cells.Item(169, 10).Value = 0
$ws.Cells.Item(169, 11).Value = 0
$ws.Cells.Item(169, 12).Value = 0

$ws.Cells.Item(170, 1).Value = "Du"
$ws.Cells.Item(170, 3).Value = 4
$ws.Cells.Item(170, 4).Value = 2
$ws.Cells.Item(170, 5).Value = 2
$ws.Cells.Item(170, 6).Value = 0
$ws.Cells.Item(170, 7).Value = 1
$ws.Cells.Item(170, 8).Value = 1
$ws.Cells.Item(170, 9).Value = 0
$ws.Cells.Item(170, 10).Value = 0
$ws.Cells.Item(170, 11).Value = 0
$ws.Cells.Item(170, 12).Value = 0

$ws.Cells.Item(171, 1).Value = "Ismael"
$ws.Cells.Item(171, 3).Value = 4
$ws.Cells.Item(171, 4).Value = 2
$ws.Cells.Item(171, 5).Value = 2
$ws.Cells.Item(171, 6).Value = 2
$ws.Cells.Item(171, 7).Value = 1
$ws.Cells.Item(171, 8).Value = 1
$ws.Cells.Item(171, 9).Value = 0
$ws.Cells.Item(171, 10).Value = 1
$ws.Cells.Item(171, 11).Value = 0
$ws.Cells.Item(171, 12).Value = 0

$ws.Cells.Item(172, 1).Value = "Athos"
$ws.Cells.Item(172, 3).Value = 2
$ws.Cells.Item(172, 4).Value = 3
$ws.Cells.Item(172, 5).Value = 2
$ws.Cells.Item(172, 6).Value = 1
$ws.Cells.Item(172, 7).Value = 1
$ws.Cells.Item(172, 8).Value = 0
$ws.Cells.Item(172, 9).Value = 0
$ws.Cells.Item(172, 10).Value = 0
$ws.Cells.Item(172, 11).Value = 0
$ws.Cells.Item(172, 12).Value = 0

$ws.Cells.Item(173, 1).Value = "Leandrinho"
$ws.Cells.Item(173, 3).Value = 2
$ws.Cells.Item(173, 4).Value = 3
$ws.Cells.Item(173, 5).Value = 2
$ws.Cells.Item(173, 6).Value = 1
$ws.Cells.Item(173, 7).Value = 1
$ws.Cells.Item(173, 8).Value = 0
$ws.Cells.Item(173, 9).Value = 0
$ws.Cells.Item(173, 10).Value = 0
$ws.Cells.Item(173, 11).Value = 0
$ws.Cells.Item(173, 12).Value = 0

$ws.Cells.Item(174, 1).Value = "Said"
$ws.Cells.Item(174, 3).Value = 2
$ws.Cells.Item(174, 4).Value = 3
$ws.Cells.Item(174, 5).Value = 2
$ws.Cells.Item(174, 6).Value = 1
$ws.Cells.Item(174, 7).Value = 1
$ws.Cells.Item(174, 8).Value = 0
$ws.Cells.Item(174, 9).Value = 0
$ws.Cells.Item(174, 10).Value = 0
$ws.Cells.Item(174, 11).Value = 0
$ws.Cells.Item(174, 12).Value = 0

$ws.Cells.Item(175, 1).Value = "Marcos"
$ws.Cells.Item(175, 3).Value = 2
$ws.Cells.Item(175, 4).Value = 3
$ws.Cells.Item(175, 5).Value = 2
$ws.Cells.Item(175, 6).Value = 1
$ws.Cells.Item(175, 7).Value = 1
$ws.Cells.Item(175, 8).Value = 0
$ws.Cells.Item(175, 9).Value = 0
$ws.Cells.Item(175, 10).Value = 0
$ws.Cells.Item(175, 11).Value = 0
$ws.Cells.Item(175, 12).Value = 0

$ws.Cells.Item(176, 1).Value = "Guinha"
$ws.Cells.Item(176, 3).Value = 2
$ws.Cells.Item(176, 4).Value = 3
$ws.Cells.Item(176, 5).Value = 2
$ws.Cells.Item(176, 6).Value = 1
$ws.Cells.Item(176, 7).Value = 1
$ws.Cells.Item(176, 8).Value = 0
$ws.Cells.Item(176, 9).Value = 0
$ws.Cells.Item(176, 10).Value = 0
$ws.Cells.Item(176, 11).Value = 0
$ws.Cells.Item(176, 12).Value = 0

$ws.Cells.Item(177, 1).Value = "Matheus"
$ws.Cells.Item(177, 3).Value = 3
$ws.Cells.Item(177, 4).Value = 5
$ws.Cells.Item(177, 5).Value = 4
$ws.Cells.Item(177, 6).Value = 0
$ws.Cells.Item(177, 7).Value = 1
$ws.Cells.Item(177, 8).Value = 0
$ws.Cells.Item(177, 9).Value = 1
$ws.Cells.Item(177, 10).Value = 0
$ws.Cells.Item(177, 11).Value = 8
$ws.Cells.Item(177, 12).Value = 0

$ws.Cells.Item(178, 1).Value = "Alan"
$ws.Cells.Item(178, 3).Value = 4
$ws.Cells.Item(178, 4).Value = 4
$ws.Cells.Item(178, 5).Value = 3
$ws.Cells.Item(178, 6).Value = 0
$ws.Cells.Item(178, 7).Value = 1
$ws.Cells.Item(178, 8).Value = 1
$ws.Cells.Item(178, 9).Value = 0
$ws.Cells.Item(178, 10).Value = 0
$ws.Cells.Item(178, 11).Value = 5
$ws.Cells.Item(178, 12).Value = 1

$ws.Cells.Item(179, 1).Value = "Chelin"
$ws.Cells.Item(179, 3).Value = 4
$ws.Cells.Item(179, 4).Value = 4
$ws.Cells.Item(179, 5).Value = 3
$ws.Cells.Item(179, 6).Value = 0
$ws.Cells.Item(179, 7).Value = 1
$ws.Cells.Item(179, 8).Value = 0
$ws.Cells.Item(179, 9).Value = 0
$ws.Cells.Item(179, 10).Value = 0
$ws.Cells.Item(179, 11).Value = 9
$ws.Cells.Item(179, 12).Value = 0

# Move the frozen-pane viewport down near the newly added rows and select C179,
# matching where the author left off editing.
$win = $excel.ActiveWindow
$win.ScrollRow = 167
$win.ScrollColumn = 1
$ws.Range("C179").Select() | Out-Null
